$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$shp = $s.Shapes.Item(2)

# Helper: always re-fetch a fresh TextRange for the whole shape text
function Get-Full {
    return $shp.TextFrame.TextRange
}

$dash = [char]0x2013

# 1) Paragraph: empty, level 2 (new paragraph after "... method")
Get-Full | Out-Null
$ins = (Get-Full).InsertAfter("`r")
$full = Get-Full
$n = $full.Paragraphs().Count
$full.Paragraphs($n).IndentLevel = 3

# 2) Paragraph: "IXmlSerializable", level 0
(Get-Full).InsertAfter("`rIXmlSerializable") | Out-Null
$full = Get-Full
$n = $full.Paragraphs().Count
$full.Paragraphs($n).IndentLevel = 1

# 3) Paragraph: "ReadXml", level 1
(Get-Full).InsertAfter("`rReadXml") | Out-Null
$full = Get-Full
$n = $full.Paragraphs().Count
$full.Paragraphs($n).IndentLevel = 2

# 4) Paragraph: "WriteXml", level 1
(Get-Full).InsertAfter("`rWriteXml") | Out-Null
$full = Get-Full
$n = $full.Paragraphs().Count
$full.Paragraphs($n).IndentLevel = 2

# 5) Paragraph: "GetSchema" + " – " + "Nop" + " – just return null", level 1
(Get-Full).InsertAfter("`rGetSchema") | Out-Null
$full = Get-Full
$n = $full.Paragraphs().Count
$full.Paragraphs($n).IndentLevel = 2
(Get-Full).InsertAfter(" $dash ") | Out-Null
(Get-Full).InsertAfter("Nop") | Out-Null
(Get-Full).InsertAfter(" $dash just return null") | Out-Null

# 6) Paragraph: empty, level 2 (final paragraph)
(Get-Full).InsertAfter("`r") | Out-Null
$full = Get-Full
$n = $full.Paragraphs().Count
$full.Paragraphs($n).IndentLevel = 3

Write-Host (Get-Full).Text
